$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 307; this shifts the existing rows 307-379 down to 308-380
# and expands the sheet dimension to A1:R380, matching the target diff.
$ws.Rows(307).Insert()

# Populate the newly inserted row 307 with the new data record.
$ws.Range("A307").Value = 3
$ws.Range("B307").Value = "Femacal de La Calera"
$ws.Range("C307").Value = "Coquimbo"
$ws.Range("D307").Value = 44754
$ws.Range("E307").Value = 5
$ws.Range("F307").Value = 100112043
$ws.Range("G307").Value = "Pepino ensalada"
$ws.Range("H307").Value = "Sin especificar"
$ws.Range("I307").Value = "Primera"
$ws.Range("J307").Value = 100
$ws.Range("K307").Value = 17000
$ws.Range("L307").Value = 18000
$ws.Range("M307").Value = 17550
$ws.Range("N307").Value = "$/caja 70 unidades"
$ws.Range("O307").Value = "Región de Arica y Parinacota"
$ws.Range("P307").Value = 251
$ws.Range("Q307").Value = 70
$ws.Range("R307").Value = "Hortaliza"

# Ensure the date cell keeps the same date number format used by the rest of column D.
$ws.Range("D307").NumberFormat = $ws.Range("D308").NumberFormat
